# Edit script implementing the commit:
# "fixed finale, all timings are extremely close, added 'little softer/louder' timings to ending"
#
# This inserts two new rows (138/139) into the finale/ending block of the
# timing sheet, fills in timing numbers (A/B columns) that had previously
# been left blank for several lyric rows, fixes up the 'finale' timing
# formula chain, and adds two new lyric/speaker lines.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fill in timing (A/B columns) for rows that previously had none.
#    These are simple "running timestamp" rows: column B holds an absolute
#    timestamp (seconds) and column A holds the delta from the previous
#    timed row via a formula.
# ---------------------------------------------------------------------
$ws.Range("B129").Value = 191.5
$ws.Range("A129").Formula = "=B129-B128"

$ws.Range("B131").Value = 196.25
$ws.Range("A131").Formula = "=B131-B129"

$ws.Range("B132").Value = 197.3
$ws.Range("A132").Formula = "=B132-B131"

# ---------------------------------------------------------------------
# 2) Insert two new rows right after row 137 ("finale") to make room for
#    the new "[Michael]" / "aahaa! Applause, applause" lines. This shifts
#    the old rows 138-145 down to 140-147, and Excel automatically keeps
#    their formulas/styles/number formats intact (adjusting relative
#    references as it goes - we explicitly rewrite the formulas below to
#    match the final, corrected reference chain).
# ---------------------------------------------------------------------
$ws.Rows("138:139").Insert()

# ---------------------------------------------------------------------
# 3) Fix up the "finale" row's timing (was a stray leftover value/formula)
# ---------------------------------------------------------------------
$ws.Range("B137").Value = 207.1
$ws.Range("A137").Formula = "=B137-B132"

# ---------------------------------------------------------------------
# 4) New row 138: "[Michael]" speaker label with its own timing
# ---------------------------------------------------------------------
$ws.Range("B138").Value = 208.9
$ws.Range("A138").Formula = "=B138-B137"

# ---------------------------------------------------------------------
# 5) New row 139: the "aahaa! Applause, applause" lyric line (no timing).
#    Row 139 keeps no A/B cells at all, so drop the blank placeholders
#    that Insert() left behind in those columns.
#    NOTE: write this shared string before "[Michael]" further up so the
#    shared-strings table ends up in the same order as the source file.
# ---------------------------------------------------------------------
$ws.Range("A139").Clear()
$ws.Range("B139").Clear()
$ws.Range("C139").Value = "aahaa! Applause, applause"
$ws.Range("C138").Value = "[Michael]"

# ---------------------------------------------------------------------
# 6) Row 140 (previously row 138): update its timing formula to point at
#    the new row 138 instead of the old row 126 reference.
# ---------------------------------------------------------------------
$ws.Range("A140").Formula = "=B140-B138"

# ---------------------------------------------------------------------
# 7) Row 143 (previously row 141): update its timing formula to point at
#    the new row 140 instead of the old row 138 reference.
# ---------------------------------------------------------------------
$ws.Range("A143").Formula = "=B143-B140"

# ---------------------------------------------------------------------
# 8) Row 145 (previously row 143): update the G-column subtraction
#    formula to reference the shifted G140:G144 cells.
# ---------------------------------------------------------------------
$ws.Range("G145").Formula = "=G140-(G141+G142+G143+G144)"

# ---------------------------------------------------------------------
# 9) Row 146 (previously row 144): update its timing formula to point at
#    the new row 143 instead of the old row 141 reference.
# ---------------------------------------------------------------------
$ws.Range("A146").Formula = "=B146-B143"

# ---------------------------------------------------------------------
# 10) Update the active selection/scroll position to match where the
#     author ended up editing (near the bottom of the sheet).
# ---------------------------------------------------------------------
$ws.Range("A141").Select()
